$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
# Remove the two trailing rows (old rows 12 and 13) so the sheet shrinks to 11 rows
$ws1.Rows.Item(12).Delete()
$ws1.Rows.Item(12).Delete()

# Rewrite rows 2-11 with the refreshed listings
# row 2
$ws1.Cells.Item(2,1).NumberFormat = "General"
$ws1.Cells.Item(2,1).Value2 = 1
$ws1.Cells.Item(2,2).NumberFormat = "@"
$ws1.Cells.Item(2,2).Value2 = '2024-11-02'
$ws1.Cells.Item(2,3).NumberFormat = "@"
$ws1.Cells.Item(2,3).Value2 = '信州·上漫·ACG动漫游戏嘉年华'
$ws1.Cells.Item(2,4).NumberFormat = "@"
$ws1.Cells.Item(2,4).Value2 = '高铁经济试验区吴楚大道与凤凰东大道交叉口 饶派沉浸式街区(B1)'
$ws1.Cells.Item(2,5).NumberFormat = "@"
$ws1.Cells.Item(2,5).Value2 = '2024.11.02 10:00-11.02 17:30'
$ws1.Cells.Item(2,6).NumberFormat = "General"
$ws1.Cells.Item(2,6).Value2 = 53
$ws1.Cells.Item(2,7).NumberFormat = "General"
$ws1.Cells.Item(2,7).Value2 = 39.9
$ws1.Cells.Item(2,8).NumberFormat = "@"
$ws1.Cells.Item(2,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=93223'
$ws1.Cells.Item(2,9).NumberFormat = "@"
$ws1.Cells.Item(2,9).Value2 = '//i0.hdslb.com/bfs/openplatform/202410/5BN1Xhzp1728294588740.jpeg'

# row 3
$ws1.Cells.Item(3,1).NumberFormat = "General"
$ws1.Cells.Item(3,1).Value2 = 2
$ws1.Cells.Item(3,2).NumberFormat = "@"
$ws1.Cells.Item(3,2).Value2 = '2024-11-02'
$ws1.Cells.Item(3,3).NumberFormat = "@"
$ws1.Cells.Item(3,3).Value2 = '南昌·花绒万兽秋镜派对'
$ws1.Cells.Item(3,4).NumberFormat = "@"
$ws1.Cells.Item(3,4).Value2 = '双港西大街899号 旭辉Cmall(南昌店)'
$ws1.Cells.Item(3,5).NumberFormat = "@"
$ws1.Cells.Item(3,5).Value2 = '2024.11.02 10:00-11.03 21:30'
$ws1.Cells.Item(3,6).NumberFormat = "General"
$ws1.Cells.Item(3,6).Value2 = 64
$ws1.Cells.Item(3,7).NumberFormat = "General"
$ws1.Cells.Item(3,7).Value2 = 168
$ws1.Cells.Item(3,8).NumberFormat = "@"
$ws1.Cells.Item(3,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=92859'
$ws1.Cells.Item(3,9).NumberFormat = "@"
$ws1.Cells.Item(3,9).Value2 = '//i2.hdslb.com/bfs/openplatform/202409/7hJL2m3F1727175584690.jpeg'

# row 4
$ws1.Cells.Item(4,1).NumberFormat = "General"
$ws1.Cells.Item(4,1).Value2 = 3
$ws1.Cells.Item(4,2).NumberFormat = "@"
$ws1.Cells.Item(4,2).Value2 = '2024-11-03'
$ws1.Cells.Item(4,3).NumberFormat = "@"
$ws1.Cells.Item(4,3).Value2 = '南昌·鸢歌弦 代号鸢同人only'
$ws1.Cells.Item(4,4).NumberFormat = "@"
$ws1.Cells.Item(4,4).Value2 = '南龙蟠街万达茂酒吧街1楼1010室在地图中查看 洛斯百年宴会艺术酒店(红谷滩旗舰店)'
$ws1.Cells.Item(4,5).NumberFormat = "@"
$ws1.Cells.Item(4,5).Value2 = '2024.11.03 10:00-11.03 20:00'
$ws1.Cells.Item(4,6).NumberFormat = "General"
$ws1.Cells.Item(4,6).Value2 = 147
$ws1.Cells.Item(4,7).NumberFormat = "General"
$ws1.Cells.Item(4,7).Value2 = 89
$ws1.Cells.Item(4,8).NumberFormat = "@"
$ws1.Cells.Item(4,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=93248'
$ws1.Cells.Item(4,9).NumberFormat = "@"
$ws1.Cells.Item(4,9).Value2 = '//i1.hdslb.com/bfs/openplatform/202409/bv3PWQhU1727247997062.jpeg'

# row 5
$ws1.Cells.Item(5,1).NumberFormat = "General"
$ws1.Cells.Item(5,1).Value2 = 4
$ws1.Cells.Item(5,2).NumberFormat = "@"
$ws1.Cells.Item(5,2).Value2 = '2024-11-16'
$ws1.Cells.Item(5,3).NumberFormat = "@"
$ws1.Cells.Item(5,3).Value2 = '上饶·星河城市动漫文化节'
$ws1.Cells.Item(5,4).NumberFormat = "@"
$ws1.Cells.Item(5,4).Value2 = '春江北大道时光PARK内 博悦宴会艺术中心'
$ws1.Cells.Item(5,5).NumberFormat = "@"
$ws1.Cells.Item(5,5).Value2 = '2024.11.16 10:00-11.16 17:00'
$ws1.Cells.Item(5,6).NumberFormat = "General"
$ws1.Cells.Item(5,6).Value2 = 351
$ws1.Cells.Item(5,7).NumberFormat = "General"
$ws1.Cells.Item(5,7).Value2 = 55
$ws1.Cells.Item(5,8).NumberFormat = "@"
$ws1.Cells.Item(5,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=92572'
$ws1.Cells.Item(5,9).NumberFormat = "@"
$ws1.Cells.Item(5,9).Value2 = '//i2.hdslb.com/bfs/openplatform/202409/xp4jNVRG1727165677359.jpeg'

# row 6
$ws1.Cells.Item(6,1).NumberFormat = "General"
$ws1.Cells.Item(6,1).Value2 = 5
$ws1.Cells.Item(6,2).NumberFormat = "@"
$ws1.Cells.Item(6,2).Value2 = '2024-11-16'
$ws1.Cells.Item(6,3).NumberFormat = "@"
$ws1.Cells.Item(6,3).Value2 = '南昌·CM04动漫游戏博览会'
$ws1.Cells.Item(6,4).NumberFormat = "@"
$ws1.Cells.Item(6,4).Value2 = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws1.Cells.Item(6,5).NumberFormat = "@"
$ws1.Cells.Item(6,5).Value2 = '2024.11.16 09:00-11.17 17:00'
$ws1.Cells.Item(6,6).NumberFormat = "General"
$ws1.Cells.Item(6,6).Value2 = 5069
$ws1.Cells.Item(6,7).NumberFormat = "General"
$ws1.Cells.Item(6,7).Value2 = 65
$ws1.Cells.Item(6,8).NumberFormat = "@"
$ws1.Cells.Item(6,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=92378'
$ws1.Cells.Item(6,9).NumberFormat = "@"
$ws1.Cells.Item(6,9).Value2 = '//i2.hdslb.com/bfs/openplatform/202409/N57Jfogr1725381095803.jpeg'

# row 7
$ws1.Cells.Item(7,1).NumberFormat = "General"
$ws1.Cells.Item(7,1).Value2 = 6
$ws1.Cells.Item(7,2).NumberFormat = "@"
$ws1.Cells.Item(7,2).Value2 = '2024-11-30'
$ws1.Cells.Item(7,3).NumberFormat = "@"
$ws1.Cells.Item(7,3).Value2 = '南昌·岁酉山河·炎国明日方舟同人ONLY'
$ws1.Cells.Item(7,4).NumberFormat = "@"
$ws1.Cells.Item(7,4).Value2 = '民德路411号 东方豪景花园酒店(民德路店)'
$ws1.Cells.Item(7,5).NumberFormat = "@"
$ws1.Cells.Item(7,5).Value2 = '2024.11.30 09:00-11.30 17:00'
$ws1.Cells.Item(7,6).NumberFormat = "General"
$ws1.Cells.Item(7,6).Value2 = 109
$ws1.Cells.Item(7,7).NumberFormat = "General"
$ws1.Cells.Item(7,7).Value2 = 68
$ws1.Cells.Item(7,8).NumberFormat = "@"
$ws1.Cells.Item(7,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=93050'
$ws1.Cells.Item(7,9).NumberFormat = "@"
$ws1.Cells.Item(7,9).Value2 = '//i2.hdslb.com/bfs/openplatform/202409/IBvdHJ1G1726720682507.png'

# row 8
$ws1.Cells.Item(8,1).NumberFormat = "General"
$ws1.Cells.Item(8,1).Value2 = 7
$ws1.Cells.Item(8,2).NumberFormat = "@"
$ws1.Cells.Item(8,2).Value2 = '2024-12-07'
$ws1.Cells.Item(8,3).NumberFormat = "@"
$ws1.Cells.Item(8,3).Value2 = '南昌·云芽动漫音乐嘉年华'
$ws1.Cells.Item(8,4).NumberFormat = "@"
$ws1.Cells.Item(8,4).Value2 = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws1.Cells.Item(8,5).NumberFormat = "@"
$ws1.Cells.Item(8,5).Value2 = '2024.12.07 09:00-12.08 18:00'
$ws1.Cells.Item(8,6).NumberFormat = "General"
$ws1.Cells.Item(8,6).Value2 = 5274
$ws1.Cells.Item(8,7).NumberFormat = "General"
$ws1.Cells.Item(8,7).Value2 = 69
$ws1.Cells.Item(8,8).NumberFormat = "@"
$ws1.Cells.Item(8,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=92144'
$ws1.Cells.Item(8,9).NumberFormat = "@"
$ws1.Cells.Item(8,9).Value2 = '//i0.hdslb.com/bfs/openplatform/202409/2DwZA4qv1725706772865.png'

# row 9
$ws1.Cells.Item(9,1).NumberFormat = "General"
$ws1.Cells.Item(9,1).Value2 = 8
$ws1.Cells.Item(9,2).NumberFormat = "@"
$ws1.Cells.Item(9,2).Value2 = '2024-12-08'
$ws1.Cells.Item(9,3).NumberFormat = "@"
$ws1.Cells.Item(9,3).Value2 = '南昌·云芽动漫音乐嘉年华·封茗囧菌内场票'
$ws1.Cells.Item(9,4).NumberFormat = "@"
$ws1.Cells.Item(9,4).Value2 = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws1.Cells.Item(9,5).NumberFormat = "@"
$ws1.Cells.Item(9,5).Value2 = '2024.12.08 09:30-12.08 17:30'
$ws1.Cells.Item(9,6).NumberFormat = "General"
$ws1.Cells.Item(9,6).Value2 = 607
$ws1.Cells.Item(9,7).NumberFormat = "General"
$ws1.Cells.Item(9,7).Value2 = 128
$ws1.Cells.Item(9,8).NumberFormat = "@"
$ws1.Cells.Item(9,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=92134'
$ws1.Cells.Item(9,9).NumberFormat = "@"
$ws1.Cells.Item(9,9).Value2 = '//i0.hdslb.com/bfs/openplatform/202409/eeFHJb3W1725328994111.jpeg'

# row 10
$ws1.Cells.Item(10,1).NumberFormat = "General"
$ws1.Cells.Item(10,1).Value2 = 9
$ws1.Cells.Item(10,2).NumberFormat = "@"
$ws1.Cells.Item(10,2).Value2 = '2025-01-01'
$ws1.Cells.Item(10,3).NumberFormat = "@"
$ws1.Cells.Item(10,3).Value2 = '南昌·萌卡动漫展'
$ws1.Cells.Item(10,4).NumberFormat = "@"
$ws1.Cells.Item(10,4).Value2 = '八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆'
$ws1.Cells.Item(10,5).NumberFormat = "@"
$ws1.Cells.Item(10,5).Value2 = '2025.01.01 09:00-01.03 17:00'
$ws1.Cells.Item(10,6).NumberFormat = "General"
$ws1.Cells.Item(10,6).Value2 = 1336
$ws1.Cells.Item(10,7).NumberFormat = "General"
$ws1.Cells.Item(10,7).Value2 = 65
$ws1.Cells.Item(10,8).NumberFormat = "@"
$ws1.Cells.Item(10,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=93031'
$ws1.Cells.Item(10,9).NumberFormat = "@"
$ws1.Cells.Item(10,9).Value2 = '//i2.hdslb.com/bfs/openplatform/202409/HTlK8fN21727112669248.jpeg'

# row 11
$ws1.Cells.Item(11,1).NumberFormat = "General"
$ws1.Cells.Item(11,1).Value2 = 10
$ws1.Cells.Item(11,2).NumberFormat = "@"
$ws1.Cells.Item(11,2).Value2 = '2025-02-14'
$ws1.Cells.Item(11,3).NumberFormat = "@"
$ws1.Cells.Item(11,3).Value2 = '九江·第二届异次元动漫嘉年华'
$ws1.Cells.Item(11,4).NumberFormat = "@"
$ws1.Cells.Item(11,4).Value2 = '长虹西大道兴城广场99号 九江半岛宾馆'
$ws1.Cells.Item(11,5).NumberFormat = "@"
$ws1.Cells.Item(11,5).Value2 = '2025.02.14 09:30-02.14 17:30'
$ws1.Cells.Item(11,6).NumberFormat = "General"
$ws1.Cells.Item(11,6).Value2 = 99
$ws1.Cells.Item(11,7).NumberFormat = "General"
$ws1.Cells.Item(11,7).Value2 = 39.8
$ws1.Cells.Item(11,8).NumberFormat = "@"
$ws1.Cells.Item(11,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=93138'
$ws1.Cells.Item(11,9).NumberFormat = "@"
$ws1.Cells.Item(11,9).Value2 = '//i1.hdslb.com/bfs/openplatform/202409/YBlAWRDD1727019019550.jpeg'

# ---- Sheet: 全部类型 (all types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
# Remove the trailing row (old row 14) so the sheet shrinks to 13 rows
$ws4.Rows.Item(14).Delete()

# Rewrite rows 3-13 with the refreshed listings (row 2 is unchanged)
# row 3
$ws4.Cells.Item(3,1).NumberFormat = "General"
$ws4.Cells.Item(3,1).Value2 = 2
$ws4.Cells.Item(3,2).NumberFormat = "@"
$ws4.Cells.Item(3,2).Value2 = '2024-11-02'
$ws4.Cells.Item(3,3).NumberFormat = "@"
$ws4.Cells.Item(3,3).Value2 = '信州·上漫·ACG动漫游戏嘉年华'
$ws4.Cells.Item(3,4).NumberFormat = "@"
$ws4.Cells.Item(3,4).Value2 = '高铁经济试验区吴楚大道与凤凰东大道交叉口 饶派沉浸式街区(B1)'
$ws4.Cells.Item(3,5).NumberFormat = "@"
$ws4.Cells.Item(3,5).Value2 = '2024.11.02 10:00-11.02 17:30'
$ws4.Cells.Item(3,6).NumberFormat = "General"
$ws4.Cells.Item(3,6).Value2 = 53
$ws4.Cells.Item(3,7).NumberFormat = "General"
$ws4.Cells.Item(3,7).Value2 = 39.9
$ws4.Cells.Item(3,8).NumberFormat = "@"
$ws4.Cells.Item(3,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=93223'
$ws4.Cells.Item(3,9).NumberFormat = "@"
$ws4.Cells.Item(3,9).Value2 = '//i0.hdslb.com/bfs/openplatform/202410/5BN1Xhzp1728294588740.jpeg'

# row 4
$ws4.Cells.Item(4,1).NumberFormat = "General"
$ws4.Cells.Item(4,1).Value2 = 3
$ws4.Cells.Item(4,2).NumberFormat = "@"
$ws4.Cells.Item(4,2).Value2 = '2024-11-02'
$ws4.Cells.Item(4,3).NumberFormat = "@"
$ws4.Cells.Item(4,3).Value2 = '南昌·花绒万兽秋镜派对'
$ws4.Cells.Item(4,4).NumberFormat = "@"
$ws4.Cells.Item(4,4).Value2 = '双港西大街899号 旭辉Cmall(南昌店)'
$ws4.Cells.Item(4,5).NumberFormat = "@"
$ws4.Cells.Item(4,5).Value2 = '2024.11.02 10:00-11.03 21:30'
$ws4.Cells.Item(4,6).NumberFormat = "General"
$ws4.Cells.Item(4,6).Value2 = 64
$ws4.Cells.Item(4,7).NumberFormat = "General"
$ws4.Cells.Item(4,7).Value2 = 168
$ws4.Cells.Item(4,8).NumberFormat = "@"
$ws4.Cells.Item(4,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=92859'
$ws4.Cells.Item(4,9).NumberFormat = "@"
$ws4.Cells.Item(4,9).Value2 = '//i2.hdslb.com/bfs/openplatform/202409/7hJL2m3F1727175584690.jpeg'

# row 5
$ws4.Cells.Item(5,1).NumberFormat = "General"
$ws4.Cells.Item(5,1).Value2 = 4
$ws4.Cells.Item(5,2).NumberFormat = "@"
$ws4.Cells.Item(5,2).Value2 = '2024-11-03'
$ws4.Cells.Item(5,3).NumberFormat = "@"
$ws4.Cells.Item(5,3).Value2 = '南昌·鸢歌弦 代号鸢同人only'
$ws4.Cells.Item(5,4).NumberFormat = "@"
$ws4.Cells.Item(5,4).Value2 = '南龙蟠街万达茂酒吧街1楼1010室在地图中查看 洛斯百年宴会艺术酒店(红谷滩旗舰店)'
$ws4.Cells.Item(5,5).NumberFormat = "@"
$ws4.Cells.Item(5,5).Value2 = '2024.11.03 10:00-11.03 20:00'
$ws4.Cells.Item(5,6).NumberFormat = "General"
$ws4.Cells.Item(5,6).Value2 = 147
$ws4.Cells.Item(5,7).NumberFormat = "General"
$ws4.Cells.Item(5,7).Value2 = 89
$ws4.Cells.Item(5,8).NumberFormat = "@"
$ws4.Cells.Item(5,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=93248'
$ws4.Cells.Item(5,9).NumberFormat = "@"
$ws4.Cells.Item(5,9).Value2 = '//i1.hdslb.com/bfs/openplatform/202409/bv3PWQhU1727247997062.jpeg'

# row 6
$ws4.Cells.Item(6,1).NumberFormat = "General"
$ws4.Cells.Item(6,1).Value2 = 5
$ws4.Cells.Item(6,2).NumberFormat = "@"
$ws4.Cells.Item(6,2).Value2 = '2024-11-06'
$ws4.Cells.Item(6,3).NumberFormat = "@"
$ws4.Cells.Item(6,3).Value2 = '南昌·松井祐贵 2024《阳光之旅》指弹吉他音乐会'
$ws4.Cells.Item(6,4).NumberFormat = "@"
$ws4.Cells.Item(6,4).Value2 = '上海路543号520Park文创公园21号01区域 瓦肆VAS NANCHANG'
$ws4.Cells.Item(6,5).NumberFormat = "@"
$ws4.Cells.Item(6,5).Value2 = '2024.11.06 20:00-11.06 21:30'
$ws4.Cells.Item(6,6).NumberFormat = "General"
$ws4.Cells.Item(6,6).Value2 = 2
$ws4.Cells.Item(6,7).NumberFormat = "General"
$ws4.Cells.Item(6,7).Value2 = 220
$ws4.Cells.Item(6,8).NumberFormat = "@"
$ws4.Cells.Item(6,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=92765'
$ws4.Cells.Item(6,9).NumberFormat = "@"
$ws4.Cells.Item(6,9).Value2 = '//i1.hdslb.com/bfs/openplatform/202409/iUNLvHVz1727082732931.jpeg'

# row 7
$ws4.Cells.Item(7,1).NumberFormat = "General"
$ws4.Cells.Item(7,1).Value2 = 6
$ws4.Cells.Item(7,2).NumberFormat = "@"
$ws4.Cells.Item(7,2).Value2 = '2024-11-16'
$ws4.Cells.Item(7,3).NumberFormat = "@"
$ws4.Cells.Item(7,3).Value2 = '上饶·星河城市动漫文化节'
$ws4.Cells.Item(7,4).NumberFormat = "@"
$ws4.Cells.Item(7,4).Value2 = '春江北大道时光PARK内 博悦宴会艺术中心'
$ws4.Cells.Item(7,5).NumberFormat = "@"
$ws4.Cells.Item(7,5).Value2 = '2024.11.16 10:00-11.16 17:00'
$ws4.Cells.Item(7,6).NumberFormat = "General"
$ws4.Cells.Item(7,6).Value2 = 351
$ws4.Cells.Item(7,7).NumberFormat = "General"
$ws4.Cells.Item(7,7).Value2 = 55
$ws4.Cells.Item(7,8).NumberFormat = "@"
$ws4.Cells.Item(7,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=92572'
$ws4.Cells.Item(7,9).NumberFormat = "@"
$ws4.Cells.Item(7,9).Value2 = '//i2.hdslb.com/bfs/openplatform/202409/xp4jNVRG1727165677359.jpeg'

# row 8
$ws4.Cells.Item(8,1).NumberFormat = "General"
$ws4.Cells.Item(8,1).Value2 = 7
$ws4.Cells.Item(8,2).NumberFormat = "@"
$ws4.Cells.Item(8,2).Value2 = '2024-11-16'
$ws4.Cells.Item(8,3).NumberFormat = "@"
$ws4.Cells.Item(8,3).Value2 = '南昌·CM04动漫游戏博览会'
$ws4.Cells.Item(8,4).NumberFormat = "@"
$ws4.Cells.Item(8,4).Value2 = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws4.Cells.Item(8,5).NumberFormat = "@"
$ws4.Cells.Item(8,5).Value2 = '2024.11.16 09:00-11.17 17:00'
$ws4.Cells.Item(8,6).NumberFormat = "General"
$ws4.Cells.Item(8,6).Value2 = 5069
$ws4.Cells.Item(8,7).NumberFormat = "General"
$ws4.Cells.Item(8,7).Value2 = 65
$ws4.Cells.Item(8,8).NumberFormat = "@"
$ws4.Cells.Item(8,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=92378'
$ws4.Cells.Item(8,9).NumberFormat = "@"
$ws4.Cells.Item(8,9).Value2 = '//i2.hdslb.com/bfs/openplatform/202409/N57Jfogr1725381095803.jpeg'

# row 9
$ws4.Cells.Item(9,1).NumberFormat = "General"
$ws4.Cells.Item(9,1).Value2 = 8
$ws4.Cells.Item(9,2).NumberFormat = "@"
$ws4.Cells.Item(9,2).Value2 = '2024-11-30'
$ws4.Cells.Item(9,3).NumberFormat = "@"
$ws4.Cells.Item(9,3).Value2 = '南昌·岁酉山河·炎国明日方舟同人ONLY'
$ws4.Cells.Item(9,4).NumberFormat = "@"
$ws4.Cells.Item(9,4).Value2 = '民德路411号 东方豪景花园酒店(民德路店)'
$ws4.Cells.Item(9,5).NumberFormat = "@"
$ws4.Cells.Item(9,5).Value2 = '2024.11.30 09:00-11.30 17:00'
$ws4.Cells.Item(9,6).NumberFormat = "General"
$ws4.Cells.Item(9,6).Value2 = 109
$ws4.Cells.Item(9,7).NumberFormat = "General"
$ws4.Cells.Item(9,7).Value2 = 68
$ws4.Cells.Item(9,8).NumberFormat = "@"
$ws4.Cells.Item(9,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=93050'
$ws4.Cells.Item(9,9).NumberFormat = "@"
$ws4.Cells.Item(9,9).Value2 = '//i2.hdslb.com/bfs/openplatform/202409/IBvdHJ1G1726720682507.png'

# row 10
$ws4.Cells.Item(10,1).NumberFormat = "General"
$ws4.Cells.Item(10,1).Value2 = 9
$ws4.Cells.Item(10,2).NumberFormat = "@"
$ws4.Cells.Item(10,2).Value2 = '2024-12-07'
$ws4.Cells.Item(10,3).NumberFormat = "@"
$ws4.Cells.Item(10,3).Value2 = '南昌·云芽动漫音乐嘉年华'
$ws4.Cells.Item(10,4).NumberFormat = "@"
$ws4.Cells.Item(10,4).Value2 = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws4.Cells.Item(10,5).NumberFormat = "@"
$ws4.Cells.Item(10,5).Value2 = '2024.12.07 09:00-12.08 18:00'
$ws4.Cells.Item(10,6).NumberFormat = "General"
$ws4.Cells.Item(10,6).Value2 = 5274
$ws4.Cells.Item(10,7).NumberFormat = "General"
$ws4.Cells.Item(10,7).Value2 = 69
$ws4.Cells.Item(10,8).NumberFormat = "@"
$ws4.Cells.Item(10,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=92144'
$ws4.Cells.Item(10,9).NumberFormat = "@"
$ws4.Cells.Item(10,9).Value2 = '//i0.hdslb.com/bfs/openplatform/202409/2DwZA4qv1725706772865.png'

# row 11
$ws4.Cells.Item(11,1).NumberFormat = "General"
$ws4.Cells.Item(11,1).Value2 = 10
$ws4.Cells.Item(11,2).NumberFormat = "@"
$ws4.Cells.Item(11,2).Value2 = '2024-12-08'
$ws4.Cells.Item(11,3).NumberFormat = "@"
$ws4.Cells.Item(11,3).Value2 = '南昌·云芽动漫音乐嘉年华·封茗囧菌内场票'
$ws4.Cells.Item(11,4).NumberFormat = "@"
$ws4.Cells.Item(11,4).Value2 = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws4.Cells.Item(11,5).NumberFormat = "@"
$ws4.Cells.Item(11,5).Value2 = '2024.12.08 09:30-12.08 17:30'
$ws4.Cells.Item(11,6).NumberFormat = "General"
$ws4.Cells.Item(11,6).Value2 = 607
$ws4.Cells.Item(11,7).NumberFormat = "General"
$ws4.Cells.Item(11,7).Value2 = 128
$ws4.Cells.Item(11,8).NumberFormat = "@"
$ws4.Cells.Item(11,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=92134'
$ws4.Cells.Item(11,9).NumberFormat = "@"
$ws4.Cells.Item(11,9).Value2 = '//i0.hdslb.com/bfs/openplatform/202409/eeFHJb3W1725328994111.jpeg'

# row 12
$ws4.Cells.Item(12,1).NumberFormat = "General"
$ws4.Cells.Item(12,1).Value2 = 11
$ws4.Cells.Item(12,2).NumberFormat = "@"
$ws4.Cells.Item(12,2).Value2 = '2025-01-01'
$ws4.Cells.Item(12,3).NumberFormat = "@"
$ws4.Cells.Item(12,3).Value2 = '南昌·萌卡动漫展'
$ws4.Cells.Item(12,4).NumberFormat = "@"
$ws4.Cells.Item(12,4).Value2 = '八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆'
$ws4.Cells.Item(12,5).NumberFormat = "@"
$ws4.Cells.Item(12,5).Value2 = '2025.01.01 09:00-01.03 17:00'
$ws4.Cells.Item(12,6).NumberFormat = "General"
$ws4.Cells.Item(12,6).Value2 = 1336
$ws4.Cells.Item(12,7).NumberFormat = "General"
$ws4.Cells.Item(12,7).Value2 = 65
$ws4.Cells.Item(12,8).NumberFormat = "@"
$ws4.Cells.Item(12,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=93031'
$ws4.Cells.Item(12,9).NumberFormat = "@"
$ws4.Cells.Item(12,9).Value2 = '//i2.hdslb.com/bfs/openplatform/202409/HTlK8fN21727112669248.jpeg'

# row 13
$ws4.Cells.Item(13,1).NumberFormat = "General"
$ws4.Cells.Item(13,1).Value2 = 12
$ws4.Cells.Item(13,2).NumberFormat = "@"
$ws4.Cells.Item(13,2).Value2 = '2025-02-14'
$ws4.Cells.Item(13,3).NumberFormat = "@"
$ws4.Cells.Item(13,3).Value2 = '九江·第二届异次元动漫嘉年华'
$ws4.Cells.Item(13,4).NumberFormat = "@"
$ws4.Cells.Item(13,4).Value2 = '长虹西大道兴城广场99号 九江半岛宾馆'
$ws4.Cells.Item(13,5).NumberFormat = "@"
$ws4.Cells.Item(13,5).Value2 = '2025.02.14 09:30-02.14 17:30'
$ws4.Cells.Item(13,6).NumberFormat = "General"
$ws4.Cells.Item(13,6).Value2 = 99
$ws4.Cells.Item(13,7).NumberFormat = "General"
$ws4.Cells.Item(13,7).Value2 = 39.8
$ws4.Cells.Item(13,8).NumberFormat = "@"
$ws4.Cells.Item(13,8).Value2 = 'https://show.bilibili.com/platform/detail.html?id=93138'
$ws4.Cells.Item(13,9).NumberFormat = "@"
$ws4.Cells.Item(13,9).Value2 = '//i1.hdslb.com/bfs/openplatform/202409/YBlAWRDD1727019019550.jpeg'
